# ------------------------------------------------------------------
# Rebuild the weekly-program grid: the 7 weekly columns (A:G) are
# reshuffled into their correct chronological week order, the week-3
# "song/prayer" row is trimmed down to just the song title, and the
# congregation-meeting block (rows 11-18) is re-extracted in fixed
# length blocks, which pushes some content down into two brand new
# rows (17 and 18).
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 2).Value = "8-14 DE ENERO"
$ws.Cells.Item(1, 4).Value = "22-28 DE ENERO"
$ws.Cells.Item(1, 5).Value = "5-11 DE FEBRERO"
$ws.Cells.Item(1, 6).Value = "12-18 DE FEBRERO"
$ws.Cells.Item(1, 7).Value = "19-25 DE FEBRERO"
# Row 2
$ws.Cells.Item(2, 2).Value = "JOB 34,35"
$ws.Cells.Item(2, 4).Value = "JOB 38,39"
$ws.Cells.Item(2, 5).Value = "SALMOS 1-4"
$ws.Cells.Item(2, 6).Value = "SALMOS 5-7"
$ws.Cells.Item(2, 7).Value = "SALMOS 8-10"
# Row 3
$ws.Cells.Item(3, 1).Value = "Canción 102"
$ws.Cells.Item(3, 2).Value = "Canción 30"
$ws.Cells.Item(3, 3).Value = "Canción 147"
$ws.Cells.Item(3, 4).Value = "Canción 11"
$ws.Cells.Item(3, 5).Value = "Canción 150"
$ws.Cells.Item(3, 6).Value = "Canción 118"
$ws.Cells.Item(3, 7).Value = "Canción 2"
# Row 5
$ws.Cells.Item(5, 2).Value = "1. ¿Le parece que la vida es injusta?"
$ws.Cells.Item(5, 4).Value = "1. ¿Dedica tiempo a observar la creación?"
$ws.Cells.Item(5, 5).Value = "1. Póngase de parte del Reino de Dios"
$ws.Cells.Item(5, 6).Value = "1. Sea leal a Jehová sin importar lo que hagan los demás"
$ws.Cells.Item(5, 7).Value = "1. “Te alabaré, oh, Jehová”"
# Row 9
$ws.Cells.Item(9, 5).Value = "4. Naturalidad: Lo que hizo Felipe"
$ws.Cells.Item(9, 6).Value = "4. Empiece conversaciones"
# Row 10
$ws.Cells.Item(10, 5).Value = "5. Naturalidad: Imite a Felipe"
$ws.Cells.Item(10, 6).Value = "5. Empiece conversaciones"
$ws.Cells.Item(10, 7).Value = "5. Haga revisitas"
# Row 11
$ws.Cells.Item(11, 1).Value = ""
$ws.Cells.Item(11, 2).Value = "6. Haga discípulos"
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = "6. Haga revisitas"
$ws.Cells.Item(11, 7).Value = "6. Discurso"
# Row 12
$ws.Cells.Item(12, 1).Value = ""
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = "7. Explique sus creencias"
$ws.Cells.Item(12, 7).Value = ""
# Row 13
$ws.Cells.Item(13, 1).Value = "NUESTRA VIDA CRISTIANA"
$ws.Cells.Item(13, 3).Value = "NUESTRA VIDA CRISTIANA"
$ws.Cells.Item(13, 4).Value = "NUESTRA VIDA CRISTIANA"
$ws.Cells.Item(13, 5).Value = "NUESTRA VIDA CRISTIANA"
$ws.Cells.Item(13, 6).Value = "NUESTRA VIDA CRISTIANA"
$ws.Cells.Item(13, 7).Value = "NUESTRA VIDA CRISTIANA"
# Row 14
$ws.Cells.Item(14, 1).Value = "Canción 116"
$ws.Cells.Item(14, 2).Value = "Canción 58"
$ws.Cells.Item(14, 3).Value = "Canción 49"
$ws.Cells.Item(14, 4).Value = "Canción 111"
$ws.Cells.Item(14, 5).Value = "Canción 32"
$ws.Cells.Item(14, 6).Value = "Canción 99"
$ws.Cells.Item(14, 7).Value = "Canción 10"
# Row 15
$ws.Cells.Item(15, 1).Value = "6. Necesidades de la congregación"
$ws.Cells.Item(15, 2).Value = "7. ¿“Predica la palabra” informalmente con entusiasmo?"
$ws.Cells.Item(15, 3).Value = "7. Esté preparado por si necesita tratamiento médico o una intervención quirúrgica"
$ws.Cells.Item(15, 4).Value = "7. Meditar en la creación nos ayuda a no perder de vista el cuadro completo"
$ws.Cells.Item(15, 5).Value = "6. Necesidades de la congregación"
$ws.Cells.Item(15, 6).Value = "8. Informe de servicio anual"
$ws.Cells.Item(15, 7).Value = "7. Cómo predicar informalmente de forma natural"
# Row 16
$ws.Cells.Item(16, 1).Value = "7. Estudio bíblico de la congregación"
$ws.Cells.Item(16, 2).Value = "8. Estudio bíblico de la congregación"
$ws.Cells.Item(16, 3).Value = "8. Estudio bíblico de la congregación"
$ws.Cells.Item(16, 4).Value = "8. Estudio bíblico de la congregación"
$ws.Cells.Item(16, 5).Value = "7. Estudio bíblico de la congregación"
$ws.Cells.Item(16, 6).Value = "9. Estudio bíblico de la congregación"
$ws.Cells.Item(16, 7).Value = "8. Necesidades de la congregación"
# Row 17
$ws.Cells.Item(17, 1).Value = "Canción 54"
$ws.Cells.Item(17, 2).Value = "Canción 138"
$ws.Cells.Item(17, 3).Value = "Canción 67"
$ws.Cells.Item(17, 4).Value = "Canción 54"
$ws.Cells.Item(17, 5).Value = "Canción 61"
$ws.Cells.Item(17, 6).Value = "Canción 83"
$ws.Cells.Item(17, 7).Value = "9. Estudio bíblico de la congregación"
# Row 18
$ws.Cells.Item(18, 1).Value = ""
$ws.Cells.Item(18, 2).Value = ""
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 6).Value = ""
$ws.Cells.Item(18, 7).Value = "Canción 65"
